$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from an existing header cell (e.g. H1) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-set values after paste (paste special formats only, but ensure text stays correct)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF)
$dataI = @(6, 3, 4, 3, 8, 6, 6, 1, 6, 1, 1, 2, 7, 6, 7)
$dataJ = @(9, 4, 9, 6, 8, 7, 8, 6, 6, 3, 4, 5, 7, 8, 8)

for ($i = 0; $i -lt 15; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
